$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 49975
$ws.Range("J75").Value = 49975
$ws.Range("L75").Value = 49975
$ws.Range("N75").Value = -51847
$ws.Range("H78").Value = 49975
$ws.Range("J78").Value = 49975
$ws.Range("L78").Value = 149925
$ws.Range("N78").Value = -159285
$ws.Range("H92").Value = 1764.7693
$ws.Range("I92").Value = 473.66666
$ws.Range("K92").Value = 473.66666
$ws.Range("M92").Value = 774.33334
$ws.Range("H112").Value = 1629
$ws.Range("J112").Value = 3000
$ws.Range("L112").Value = 9000
$ws.Range("N112").Value = -11216
$ws.Range("H113").Value = 7123.091
$ws.Range("I113").Value = 1966
$ws.Range("J113").Value = 9057
$ws.Range("K113").Value = 1966
$ws.Range("L113").Value = 9057
$ws.Range("M113").Value = 1288
$ws.Range("N113").Value = -15565
$ws.Range("H132").Value = 1266.1765
$ws.Range("I132").Value = 1202.9667
$ws.Range("J132").Value = 1740.25
$ws.Range("K132").Value = 3608.9001
$ws.Range("L132").Value = 5220.75
$ws.Range("M132").Value = -1078.9001
$ws.Range("N132").Value = -10280.75
$ws.Range("H138").Value = 3031.7878
$ws.Range("I138").Value = 1918.4584
$ws.Range("J138").Value = 3667.976
$ws.Range("K138").Value = 5755.3752
$ws.Range("L138").Value = 11003.928
$ws.Range("M138").Value = -615.3752000000004
$ws.Range("N138").Value = -21283.928
$ws.Range("H141").Value = 7750
$ws.Range("I141").Value = 2000
$ws.Range("J141").Value = 8571.429
$ws.Range("K141").Value = 6000
$ws.Range("L141").Value = 25714.287
$ws.Range("M141").Value = -820
$ws.Range("N141").Value = -36074.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4736.1396
$ws.Range("I32").Value = 4015.2683
$ws.Range("J32").Value = 19514
$ws.Range("K32").Value = 4015.2683
$ws.Range("L32").Value = 19514
$ws.Range("M32").Value = -3728.2683
$ws.Range("N32").Value = -20088
$ws.Range("H45").Value = 83338800
$ws.Range("I45").Value = 166667760
$ws.Range("J45").Value = 9835.333000000001
$ws.Range("K45").Value = 166667760
$ws.Range("L45").Value = 9835.333000000001
$ws.Range("M45").Value = -166667383
$ws.Range("N45").Value = -10589.333
$ws.Range("H74").Value = 18520960
$ws.Range("I74").Value = 27780466
$ws.Range("K74").Value = 27780466
$ws.Range("M74").Value = -27779592
$ws.Range("H77").Value = 18520960
$ws.Range("I77").Value = 27780466
$ws.Range("K77").Value = 138902330
$ws.Range("M77").Value = -138897962
$ws.Range("H110").Value = 4349.3125
$ws.Range("I110").Value = 1427.9231
$ws.Range("K110").Value = 1427.9231
$ws.Range("M110").Value = 617.0769
$ws.Range("H132").Value = 2831.889
$ws.Range("I132").Value = 2246.9546
$ws.Range("J132").Value = 5405.6
$ws.Range("K132").Value = 6740.8638
$ws.Range("L132").Value = 16216.8
$ws.Range("M132").Value = -4210.8638
$ws.Range("N132").Value = -21276.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2741.6667
$ws.Range("I20").Value = 2193.5
$ws.Range("K20").Value = 2193.5
$ws.Range("M20").Value = -1946.5
$ws.Range("H80").Value = 1188.0834
$ws.Range("J80").Value = 939.7143
$ws.Range("L80").Value = 939.7143
$ws.Range("N80").Value = -2935.7143
$ws.Range("H83").Value = 1188.0834
$ws.Range("J83").Value = 939.7143
$ws.Range("L83").Value = 4698.5715
$ws.Range("N83").Value = -14682.5715
$ws.Range("H99").Value = 1614.125
$ws.Range("I99").Value = 1614.125
$ws.Range("K99").Value = 1614.125
$ws.Range("M99").Value = -116.125
$ws.Range("H134").Value = 1773.5714
$ws.Range("I134").Value = 1091.1052
$ws.Range("K134").Value = 3273.3156
$ws.Range("M134").Value = -738.3155999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4549.3
$ws.Range("I58").Value = 2155.25
$ws.Range("K58").Value = 2155.25
$ws.Range("M58").Value = -1952.25
$ws.Range("H86").Value = 7737.25
$ws.Range("I86").Value = 6266.909
$ws.Range("K86").Value = 6266.909
$ws.Range("M86").Value = -5143.909
$ws.Range("H89").Value = 7737.25
$ws.Range("I89").Value = 6266.909
$ws.Range("K89").Value = 31334.545
$ws.Range("M89").Value = -25718.545
$ws.Range("H107").Value = 803.30304
$ws.Range("I107").Value = 680.3182
$ws.Range("J107").Value = 1049.2727
$ws.Range("K107").Value = 680.3182
$ws.Range("L107").Value = 1049.2727
$ws.Range("M107").Value = 1239.6818
$ws.Range("N107").Value = -4889.2727
$ws.Range("H132").Value = 10650.8
$ws.Range("I132").Value = 10528.306
$ws.Range("J132").Value = 11753.25
$ws.Range("K132").Value = 31584.918
$ws.Range("L132").Value = 35259.75
$ws.Range("M132").Value = -29054.918
$ws.Range("N132").Value = -40319.75
$ws.Range("H134").Value = 3271.8823
$ws.Range("I134").Value = 2157
$ws.Range("K134").Value = 6471
$ws.Range("M134").Value = -3936
$ws.Range("H136").Value = 4549.3
$ws.Range("I136").Value = 2155.25
$ws.Range("K136").Value = 6465.75
$ws.Range("M136").Value = -3915.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 699.5
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H132").Value = 3469.3845
$ws.Range("J132").Value = 3639
$ws.Range("L132").Value = 32751
$ws.Range("N132").Value = -37811
$ws.Range("H137").Value = 4183.1665
$ws.Range("I137").Value = 1896.25
$ws.Range("J137").Value = 5326.625
$ws.Range("K137").Value = 5688.75
$ws.Range("L137").Value = 15979.875
$ws.Range("M137").Value = -588.75
$ws.Range("N137").Value = -26179.875
$ws.Range("H139").Value = 3140.625
$ws.Range("I139").Value = 1694.0714
$ws.Range("J139").Value = 13266.5
$ws.Range("K139").Value = 5082.2142
$ws.Range("L139").Value = 39799.5
$ws.Range("M139").Value = 57.78579999999965
$ws.Range("N139").Value = -50079.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5096.636
$ws.Range("I80").Value = 2871
$ws.Range("K80").Value = 2871
$ws.Range("M80").Value = -1873
$ws.Range("H83").Value = 5096.636
$ws.Range("I83").Value = 2871
$ws.Range("K83").Value = 14355
$ws.Range("M83").Value = -9363
$ws.Range("H102").Value = 3198.6
$ws.Range("I102").Value = 3248.5
$ws.Range("K102").Value = 3248.5
$ws.Range("M102").Value = -1626.5
$ws.Range("H122").Value = 4340.8335
$ws.Range("I122").Value = 3909.3845
$ws.Range("K122").Value = 11728.1535
$ws.Range("M122").Value = -9278.1535
$ws.Range("H132").Value = 46547.125
$ws.Range("I132").Value = 74346.5
$ws.Range("J132").Value = 7628
$ws.Range("K132").Value = 223039.5
$ws.Range("L132").Value = 22884
$ws.Range("M132").Value = -220509.5
$ws.Range("N132").Value = -27944

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3665.25
$ws.Range("J46").Value = 3999.3
$ws.Range("L46").Value = 3999.3
$ws.Range("N46").Value = -4375.3
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H93").Value = 3605.0625
$ws.Range("I93").Value = 3316.7
$ws.Range("J93").Value = 4085.6667
$ws.Range("K93").Value = 3316.7
$ws.Range("L93").Value = 4085.6667
$ws.Range("M93").Value = -2068.7
$ws.Range("N93").Value = -6581.6667
$ws.Range("H132").Value = 3277.1843
$ws.Range("I132").Value = 2040.3846
$ws.Range("J132").Value = 5956.9165
$ws.Range("K132").Value = 6121.1538
$ws.Range("L132").Value = 17870.7495
$ws.Range("M132").Value = -3591.1538
$ws.Range("N132").Value = -22930.7495
$ws.Range("H136").Value = 9313.23
$ws.Range("I136").Value = 3106.5
$ws.Range("K136").Value = 9319.5
$ws.Range("M136").Value = -6769.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 2133336.8
$ws.Range("I29").Value = 10
$ws.Range("K29").Value = 10
$ws.Range("M29").Value = 280
$ws.Range("H81").Value = 3655.5386
$ws.Range("I81").Value = 2567.5
$ws.Range("K81").Value = 5135
$ws.Range("M81").Value = -4074
$ws.Range("H84").Value = 3655.5386
$ws.Range("I84").Value = 2567.5
$ws.Range("K84").Value = 25675
$ws.Range("M84").Value = -20371
$ws.Range("H96").Value = 3964.5
$ws.Range("J96").Value = 5674.75
$ws.Range("L96").Value = 5674.75
$ws.Range("N96").Value = -8420.75
$ws.Range("H107").Value = 526
$ws.Range("I107").Value = 473.875
$ws.Range("K107").Value = 1421.625
$ws.Range("M107").Value = 498.375
$ws.Range("H113").Value = 628.2941
$ws.Range("I113").Value = 306.8889
$ws.Range("J113").Value = 989.875
$ws.Range("K113").Value = 920.6667
$ws.Range("L113").Value = 2969.625
$ws.Range("M113").Value = 1249.3333
$ws.Range("N113").Value = -7309.625
$ws.Range("H132").Value = 5123.2
$ws.Range("I132").Value = 4758.9546
$ws.Range("K132").Value = 14276.8638
$ws.Range("M132").Value = -11746.8638
$ws.Range("H135").Value = 65972.3
$ws.Range("J135").Value = 65972.3
$ws.Range("L135").Value = 65972.3
$ws.Range("N135").Value = -76112.3
$ws.Range("H136").Value = 1624
$ws.Range("I136").Value = 1624
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4872
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2322
$ws.Range("N136").ClearContents()

Write-Output "All edits applied"
